# Update the "想去人数" (column F) values for rows 2-6 on the "展览" and
# "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 6412
    $ws.Range("F3").Value = 29
    $ws.Range("F4").Value = 187
    $ws.Range("F5").Value = 1011
    $ws.Range("F6").Value = 113
}
